$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7faa9f73c820>),
                ('model',
                 AdaBoostClassifier(estimator=LGBMClassifier(boosting_type='dart',
                                                             colsample_bytree=0.7,
                                                             learning_rate=0.01,
                                                             max_depth=3,
                                                             num_leaves=2,
                                                             random_state=42,
                                                             subsample=0.7),
                                    n_estimators=5, random_state=42))])
'@

$ws.Range("B2").Value = 0.6732733932733932

$ws.Range("C2").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7faa9f69a070>, 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__subsample': 0.7, 'model__estimator__num_leaves': 2, 'model__estimator__min_child_samples': 20, 'model__estimator__max_depth': 3, 'model__estimator__learning_rate': 0.01, 'model__estimator__colsample_bytree': 0.7, 'model__estimator__class_weight': None, 'model__estimator__boosting_type': 'dart'}
'@

$ws.Range("D2").Value = 0.9612130063717511

$ws.Range("E2").Value = 0.5501962814962815

$ws.Range("F2").Value = 0.787878787878788

$ws.Range("G2").Value = 0.9570698995220637

$ws.Range("H2").Value = 0.5527142857142857

$ws.Range("I2").Value = 0.7647058823529411

$ws.Range("J2").Value = 0.9665106382978722

$ws.Range("K2").Value = 0.5683333333333332

$ws.Range("L2").Value = 0.8125

$ws.Range("M2").Value = @'
[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]
'@

$ws.Range("N2").Value = @'
[1 1 1 1 1 1 0 0 1 1 0 1 1 1 0 0 0 1 0 1 1 1 1 1]
'@

$ws.Range("A3").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7faa9f766610>),
                ('model',
                 AdaBoostClassifier(estimator=LGBMClassifier(boosting_type='dart',
                                                             colsample_bytree=0.9,
                                                             learning_rate=0.2,
                                                             max_depth=5,
                                                             num_leaves=5,
                                                             random_state=42,
                                                             subsample=0.5),
                                    n_estimators=10, random_state=42))])
'@

$ws.Range("B3").Value = 0.6192840492840492

$ws.Range("C3").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7faa9f69af40>, 'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__subsample': 0.5, 'model__estimator__num_leaves': 5, 'model__estimator__min_child_samples': 20, 'model__estimator__max_depth': 5, 'model__estimator__learning_rate': 0.2, 'model__estimator__colsample_bytree': 0.9, 'model__estimator__class_weight': None, 'model__estimator__boosting_type': 'dart'}
'@

$ws.Range("D3").Value = 0.9564773669395636

$ws.Range("E3").Value = 0.4997880896880897

$ws.Range("F3").Value = 0.6896551724137931

$ws.Range("G3").Value = 0.9588077222520017

$ws.Range("H3").Value = 0.5630944444444445

$ws.Range("I3").Value = 0.7692307692307693

$ws.Range("J3").Value = 0.9560425531914893

$ws.Range("K3").Value = 0.473

$ws.Range("M3").Value = @'
[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]
'@

$ws.Range("N3").Value = @'
[0 0 0 1 1 0 1 1 1 1 1 0 0 0 1 1 1 1 1 0 0 0 0 1]
'@

$ws.Range("A4").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 AdaBoostClassifier(estimator=LGBMClassifier(boosting_type='dart',
                                                             colsample_bytree=0.9,
                                                             learning_rate=0.01,
                                                             max_depth=1,
                                                             num_leaves=2,
                                                             random_state=42,
                                                             subsample=0.5),
                                    n_estimators=5, random_state=42))])
'@

$ws.Range("B4").Value = 0.6468220668220669

$ws.Range("C4").Value = @'
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__subsample': 0.5, 'model__estimator__num_leaves': 2, 'model__estimator__min_child_samples': 20, 'model__estimator__max_depth': 1, 'model__estimator__learning_rate': 0.01, 'model__estimator__colsample_bytree': 0.9, 'model__estimator__class_weight': None, 'model__estimator__boosting_type': 'dart'}
'@

$ws.Range("D4").Value = 0.9565216223162405

$ws.Range("E4").Value = 0.5214840159840162

$ws.Range("F4").Value = 0.6666666666666667

$ws.Range("G4").Value = 0.9563322359432306

$ws.Range("H4").Value = 0.5599880952380953

$ws.Range("I4").Value = 0.7857142857142857

$ws.Range("J4").Value = 0.9574666666666667

$ws.Range("K4").Value = 0.514

$ws.Range("L4").Value = 0.5789473684210527

$ws.Range("M4").Value = @'
[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]
'@

$ws.Range("N4").Value = @'
[0 1 1 1 0 1 0 1 1 0 0 0 0 1 1 1 0 0 1 0 1 1 1 1]
'@

$ws.Range("A5").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 AdaBoostClassifier(estimator=LGBMClassifier(colsample_bytree=0.9,
                                                             learning_rate=0.01,
                                                             max_depth=7,
                                                             min_child_samples=10,
                                                             num_leaves=10,
                                                             random_state=42,
                                                             subsample=0.5),
                                    n_estimators=10, random_state=42))])
'@

$ws.Range("B5").Value = 0.6883699633699634

$ws.Range("C5").Value = @'
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__subsample': 0.5, 'model__estimator__num_leaves': 10, 'model__estimator__min_child_samples': 10, 'model__estimator__max_depth': 7, 'model__estimator__learning_rate': 0.01, 'model__estimator__colsample_bytree': 0.9, 'model__estimator__class_weight': None, 'model__estimator__boosting_type': 'gbdt'}
'@

$ws.Range("D5").Value = 0.9670418829338966

$ws.Range("E5").Value = 0.5799967476967477

$ws.Range("F5").Value = 0.5517241379310344

$ws.Range("G5").Value = 0.9643863293515125

$ws.Range("H5").Value = 0.5779920634920636

$ws.Range("I5").Value = 0.5333333333333333

$ws.Range("J5").Value = 0.9722857142857143

$ws.Range("K5").Value = 0.6046666666666667

$ws.Range("L5").Value = 0.5714285714285714

$ws.Range("M5").Value = @'
[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]
'@

$ws.Range("N5").Value = @'
[0 0 1 0 1 1 1 1 0 1 1 0 0 1 1 0 0 1 1 0 1 1 1 1]
'@

$ws.Range("A6").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7faa9f75f400>),
                ('model',
                 AdaBoostClassifier(estimator=LGBMClassifier(class_weight='balanced',
                                                             colsample_bytree=0.9,
                                                             max_depth=5,
                                                             num_leaves=5,
                                                             random_state=42,
                                                             subsample=0.9),
                                    n_estimators=10, random_state=42))])
'@

$ws.Range("B6").Value = 0.703076923076923

$ws.Range("C6").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7faa9f43b850>, 'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__subsample': 0.9, 'model__estimator__num_leaves': 5, 'model__estimator__min_child_samples': 20, 'model__estimator__max_depth': 5, 'model__estimator__learning_rate': 0.1, 'model__estimator__colsample_bytree': 0.9, 'model__estimator__class_weight': 'balanced', 'model__estimator__boosting_type': 'gbdt'}
'@

$ws.Range("D6").Value = 0.9543629484469506

$ws.Range("E6").Value = 0.6146100788100788

$ws.Range("F6").Value = 0.6206896551724137

$ws.Range("G6").Value = 0.9576397749305066

$ws.Range("H6").Value = 0.6100984126984126

$ws.Range("I6").Value = 0.5

$ws.Range("J6").Value = 0.953923076923077

$ws.Range("K6").Value = 0.6406666666666667

$ws.Range("L6").Value = 0.8181818181818182

$ws.Range("M6").Value = @'
[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]
'@

$ws.Range("N6").Value = @'
[1 0 1 1 1 1 1 1 0 0 1 1 1 1 0 0 1 1 1 0 1 1 1 1]
'@
